# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Macroferia Regional de Talca - Membrillo) just
# before what used to be row 45, pushing the former rows 45-72 down to 48-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at position 45; Excel automatically shifts rows 45-72
# down to 48-75 and extends the used range / dimension accordingly.
$ws.Rows("45:47").Insert()

# --- New row 45 ---
$ws.Range("A45").Value = 5
$ws.Range("B45").Value = "Macroferia Regional de Talca"
$ws.Range("C45").Value = "Maule"
$ws.Range("D45").Value = 45033
$ws.Range("E45").Value = 7
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100104
$ws.Range("H45").Value = "Frutos de pepita"
$ws.Range("I45").Value = 100104003
$ws.Range("J45").Value = "Membrillo"
$ws.Range("K45").Value = "Champion"
$ws.Range("L45").Value = "Especial"
$ws.Range("M45").Value = 230
$ws.Range("N45").Value = 10000
$ws.Range("O45").Value = 10000
$ws.Range("P45").Value = 10000
$ws.Range("Q45").Value = "`$/caja 18 kilos granel"
$ws.Range("R45").Value = "Región de O'Higgins"
$ws.Range("S45").Value = 556
$ws.Range("T45").Value = 18

# --- New row 46 ---
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 45033
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100104
$ws.Range("H46").Value = "Frutos de pepita"
$ws.Range("I46").Value = 100104003
$ws.Range("J46").Value = "Membrillo"
$ws.Range("K46").Value = "Champion"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 250
$ws.Range("N46").Value = 8000
$ws.Range("O46").Value = 8000
$ws.Range("P46").Value = 8000
$ws.Range("Q46").Value = "`$/caja 18 kilos granel"
$ws.Range("R46").Value = "Región de O'Higgins"
$ws.Range("S46").Value = 444
$ws.Range("T46").Value = 18

# --- New row 47 ---
$ws.Range("A47").Value = 5
$ws.Range("B47").Value = "Macroferia Regional de Talca"
$ws.Range("C47").Value = "Maule"
$ws.Range("D47").Value = 45033
$ws.Range("E47").Value = 7
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100104
$ws.Range("H47").Value = "Frutos de pepita"
$ws.Range("I47").Value = 100104003
$ws.Range("J47").Value = "Membrillo"
$ws.Range("K47").Value = "Champion"
$ws.Range("L47").Value = "Segunda"
$ws.Range("M47").Value = 200
$ws.Range("N47").Value = 6000
$ws.Range("O47").Value = 6000
$ws.Range("P47").Value = 6000
$ws.Range("Q47").Value = "`$/caja 18 kilos granel"
$ws.Range("R47").Value = "Región de O'Higgins"
$ws.Range("S47").Value = 333
$ws.Range("T47").Value = 18
